# G1_LP_A.xlsx — "changed gfx size, enable all keys, correct LP phases (resorted)"
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# The LP-phase rows were resorted; rewrite A:F for every row whose
# (start, key1, key2, day, key3, shift) tuple moved to a new row.
$rows = @(
    @{Row=8;  A=22; B="pe"; C="lo"; D="sa"; E="fu"; F="to"},
    @{Row=9;  A=38; B="fi"; C="lo"; D="la"; E="fu"; F="se"},
    @{Row=10; A=19; B="pe"; C="lo"; D="ra"; E="fu"; F="to"},
    @{Row=11; A=34; B="fi"; C="lo"; D="fe"; E="fu"; F="se"},
    @{Row=12; A=43; B="fi"; C="lo"; D="ni"; E="fu"; F="se"},
    @{Row=13; A=24; B="pe"; C="lo"; D="ta"; E="fu"; F="to"},
    @{Row=14; A=55; B="fi"; C="lo"; D="wa"; E="fu"; F="se"},
    @{Row=15; A=47; B="fi"; C="lo"; D="ra"; E="fu"; F="se"},
    @{Row=16; A=8;  B="pe"; C="lo"; D="ki"; E="fu"; F="to"},
    @{Row=17; A=53; B="fi"; C="lo"; D="te"; E="fu"; F="se"},
    @{Row=18; A=7;  B="pe"; C="lo"; D="ka"; E="fu"; F="to"},
    @{Row=19; A=44; B="fi"; C="lo"; D="pi"; E="fu"; F="se"},
    @{Row=20; A=26; B="pe"; C="lo"; D="ti"; E="fu"; F="to"},
    @{Row=21; A=36; B="fi"; C="lo"; D="ki"; E="fu"; F="se"},
    @{Row=22; A=3;  B="pe"; C="lo"; D="di"; E="fu"; F="to"},
    @{Row=23; A=46; B="fi"; C="lo"; D="po"; E="fu"; F="se"},
    @{Row=24; A=17; B="pe"; C="lo"; D="pa"; E="fu"; F="to"},
    @{Row=25; A=35; B="fi"; C="lo"; D="ka"; E="fu"; F="se"},
    @{Row=30; A=49; B="fi"; C="lo"; D="ri"; E="fu"; F="se"},
    @{Row=31; A=25; B="pe"; C="lo"; D="te"; E="fu"; F="to"},
    @{Row=32; A=41; B="fi"; C="lo"; D="mi"; E="fu"; F="se"},
    @{Row=33; A=1;  B="pe"; C="lo"; D="bo"; E="fu"; F="to"}
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 1).Value = $r.A
    $ws.Cells.Item($r.Row, 2).Value = $r.B
    $ws.Cells.Item($r.Row, 3).Value = $r.C
    $ws.Cells.Item($r.Row, 4).Value = $r.D
    $ws.Cells.Item($r.Row, 5).Value = $r.E
    $ws.Cells.Item($r.Row, 6).Value = $r.F
}

# "changed gfx size": zoom level bumped from 62% to 90%, and the
# last-used selection moved from P18 to C29.
$excel.ActiveWindow.Zoom = 90
$null = $ws.Range("C29").Select()
